$d = $word.ActiveDocument

# The third table (Commando/Parameters/Omschrijving) has a last row
# describing the "fill" command which needs to be removed.
$table = $d.Tables.Item(3)
$lastRow = $table.Rows.Item($table.Rows.Count)
$lastRow.Delete()
